$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 44.471632
$ws.Range("H2").Value = 133.414896
$ws.Range("I2").Value = 0.2668714234083088
$ws.Range("J2").Value = 0.2668714234083088
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.17516666666667
$ws.Range("N2").Value = 48.52549999999999
$ws.Range("O2").Value = 0.2200960555056711
$ws.Range("P2").Value = 0.2200960555056711
$ws.Range("Q2").Value = 719.3360595386666
$ws.Range("R2").Value = 6474.024535847999
$ws.Range("S2").Value = 0.0587373476193526
$ws.Range("T2").Value = 0.0587373476193526

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 44.471632
$ws.Range("H3").Value = 133.414896
$ws.Range("I3").Value = 0.2668714234083088
$ws.Range("J3").Value = 0.2668714234083088
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 21.45763766666667
$ws.Range("N3").Value = 64.372913
$ws.Range("O3").Value = 0.2919748221596838
$ws.Range("P3").Value = 0.2919748221596839
$ws.Range("Q3").Value = 954.2561659013387
$ws.Range("R3").Value = 8588.305493112048
$ws.Range("S3").Value = 0.07791973638914265
$ws.Range("T3").Value = 0.07791973638914268

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 44.471632
$ws.Range("H4").Value = 133.414896
$ws.Range("I4").Value = 0.2668714234083088
$ws.Range("J4").Value = 0.2668714234083088
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.29472433333333
$ws.Range("N4").Value = 48.884173
$ws.Range("O4").Value = 0.2217228808349596
$ws.Range("P4").Value = 0.2217228808349596
$ws.Range("Q4").Value = 724.6529840934454
$ws.Range("R4").Value = 6521.876856841009
$ws.Range("S4").Value = 0.0591715008106165
$ws.Range("T4").Value = 0.05917150081061651

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 44.471632
$ws.Range("H5").Value = 133.414896
$ws.Range("I5").Value = 0.2668714234083088
$ws.Range("J5").Value = 0.2668714234083088
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.56386866666667
$ws.Range("N5").Value = 58.69160600000001
$ws.Range("O5").Value = 0.2662062414996854
$ws.Range("P5").Value = 0.2662062414996854
$ws.Range("Q5").Value = 870.0371678403308
$ws.Range("R5").Value = 7830.334510562977
$ws.Range("S5").Value = 0.07104283858919705
$ws.Range("T5").Value = 0.07104283858919706

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.05790966666667
$ws.Range("H6").Value = 111.173729
$ws.Range("I6").Value = 0.2223821491705063
$ws.Range("J6").Value = 0.2223821491705063
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 16.17516666666667
$ws.Range("N6").Value = 48.52549999999999
$ws.Range("O6").Value = 0.2200960555056711
$ws.Range("P6").Value = 0.2200960555056711
$ws.Range("Q6").Value = 599.4178651766111
$ws.Range("R6").Value = 5394.760786589499
$ws.Range("S6").Value = 0.04894543384730219
$ws.Range("T6").Value = 0.04894543384730218

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.05790966666667
$ws.Range("H7").Value = 111.173729
$ws.Range("I7").Value = 0.2223821491705063
$ws.Range("J7").Value = 0.2223821491705063
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.45763766666667
$ws.Range("N7").Value = 64.372913
$ws.Range("O7").Value = 0.2919748221596838
$ws.Range("P7").Value = 0.2919748221596839
$ws.Range("Q7").Value = 795.1751983113975
$ws.Range("R7").Value = 7156.576784802576
$ws.Range("S7").Value = 0.06492998845554686
$ws.Range("T7").Value = 0.06492998845554687

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.05790966666667
$ws.Range("H8").Value = 111.173729
$ws.Range("I8").Value = 0.2223821491705063
$ws.Range("J8").Value = 0.2223821491705063
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.29472433333333
$ws.Range("N8").Value = 48.884173
$ws.Range("O8").Value = 0.2217228808349596
$ws.Range("P8").Value = 0.2217228808349596
$ws.Range("Q8").Value = 603.8484223879019
$ws.Range("R8").Value = 5434.635801491117
$ws.Range("S8").Value = 0.04930721076035437
$ws.Range("T8").Value = 0.04930721076035437

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.05790966666667
$ws.Range("H9").Value = 111.173729
$ws.Range("I9").Value = 0.2223821491705063
$ws.Range("J9").Value = 0.2223821491705063
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.56386866666667
$ws.Range("N9").Value = 58.69160600000001
$ws.Range("O9").Value = 0.2662062414996854
$ws.Range("P9").Value = 0.2662062414996854
$ws.Range("Q9").Value = 724.9960777798639
$ws.Range("R9").Value = 6524.964700018774
$ws.Range("S9").Value = 0.05919951610730285
$ws.Range("T9").Value = 0.05919951610730286

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 39.65536166666666
$ws.Range("H10").Value = 118.966085
$ws.Range("I10").Value = 0.2379692927337279
$ws.Range("J10").Value = 0.2379692927337279
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.17516666666667
$ws.Range("N10").Value = 48.52549999999999
$ws.Range("O10").Value = 0.2200960555056711
$ws.Range("P10").Value = 0.2200960555056711
$ws.Range("Q10").Value = 641.4320841852777
$ws.Range("R10").Value = 5772.888757667499
$ws.Range("S10").Value = 0.05237610266216787
$ws.Range("T10").Value = 0.05237610266216787

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 39.65536166666666
$ws.Range("H11").Value = 118.966085
$ws.Range("I11").Value = 0.2379692927337279
$ws.Range("J11").Value = 0.2379692927337279
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 21.45763766666667
$ws.Range("N11").Value = 64.372913
$ws.Range("O11").Value = 0.2919748221596838
$ws.Range("P11").Value = 0.2919748221596839
$ws.Range("Q11").Value = 850.910382183956
$ws.Range("R11").Value = 7658.193439655604
$ws.Range("S11").Value = 0.06948104192539593
$ws.Range("T11").Value = 0.06948104192539595

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 39.65536166666666
$ws.Range("H12").Value = 118.966085
$ws.Range("I12").Value = 0.2379692927337279
$ws.Range("J12").Value = 0.2379692927337279
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 16.29472433333333
$ws.Range("N12").Value = 48.884173
$ws.Range("O12").Value = 0.2217228808349596
$ws.Range("P12").Value = 0.2217228808349596
$ws.Range("Q12").Value = 646.1731866969673
$ws.Range("R12").Value = 5815.558680272705
$ws.Range("S12").Value = 0.05276323713517996
$ws.Range("T12").Value = 0.05276323713517996

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 39.65536166666666
$ws.Range("H13").Value = 118.966085
$ws.Range("I13").Value = 0.2379692927337279
$ws.Range("J13").Value = 0.2379692927337279
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.56386866666667
$ws.Range("N13").Value = 58.69160600000001
$ws.Range("O13").Value = 0.2662062414996854
$ws.Range("P13").Value = 0.2662062414996854
$ws.Range("Q13").Value = 775.8122875758345
$ws.Range("R13").Value = 6982.31058818251
$ws.Range("S13").Value = 0.06334891101098408
$ws.Range("T13").Value = 0.0633489110109841

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 45.45576366666668
$ws.Range("H14").Value = 136.367291
$ws.Range("I14").Value = 0.272777134687457
$ws.Range("J14").Value = 0.272777134687457
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 16.17516666666667
$ws.Range("N14").Value = 48.52549999999999
$ws.Range("O14").Value = 0.2200960555056711
$ws.Range("P14").Value = 0.2200960555056711
$ws.Range("Q14").Value = 735.2545532689446
$ws.Range("R14").Value = 6617.290979420501
$ws.Range("S14").Value = 0.06003717137684847
$ws.Range("T14").Value = 0.06003717137684847

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 45.45576366666668
$ws.Range("H15").Value = 136.367291
$ws.Range("I15").Value = 0.272777134687457
$ws.Range("J15").Value = 0.272777134687457
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 21.45763766666667
$ws.Range("N15").Value = 64.372913
$ws.Range("O15").Value = 0.2919748221596838
$ws.Range("P15").Value = 0.2919748221596839
$ws.Range("Q15").Value = 975.373306620965
$ws.Range("R15").Value = 8778.359759588684
$ws.Range("S15").Value = 0.07964405538959839
$ws.Range("T15").Value = 0.0796440553895984

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 45.45576366666668
$ws.Range("H16").Value = 136.367291
$ws.Range("I16").Value = 0.272777134687457
$ws.Range("J16").Value = 0.272777134687457
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.29472433333333
$ws.Range("N16").Value = 48.884173
$ws.Range("O16").Value = 0.2217228808349596
$ws.Range("P16").Value = 0.2217228808349596
$ws.Range("Q16").Value = 740.6891383094828
$ws.Range("R16").Value = 6666.202244785345
$ws.Range("S16").Value = 0.06048093212880875
$ws.Range("T16").Value = 0.06048093212880876

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 45.45576366666668
$ws.Range("H17").Value = 136.367291
$ws.Range("I17").Value = 0.272777134687457
$ws.Range("J17").Value = 0.272777134687457
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.56386866666667
$ws.Range("N17").Value = 58.69160600000001
$ws.Range("O17").Value = 0.2662062414996854
$ws.Range("P17").Value = 0.2662062414996854
$ws.Range("Q17").Value = 889.2905905177054
$ws.Range("R17").Value = 8003.615314659349
$ws.Range("S17").Value = 0.07261497579220139
$ws.Range("T17").Value = 0.0726149757922014
